$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

# Column G ("Recorded By") holds values like "System, dnasr281@gmail.com"
# for every recorded session row. Swap the author order to
# "dnasr281@gmail.com, System" everywhere it occurs in that column.
$col = $ws.Range("G:G")
$col.Replace("System, dnasr281@gmail.com", "dnasr281@gmail.com, System")
